# Slide 1, Title placeholder: "The TODO list that doesn't freak you out"
#   -> "The TODO list that doesn't burn you out"
# (commit message: "changed 'freak' to 'burn'")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# 1) Swap the word itself: freak -> burn
$text = $tr.Text
$idx = $text.IndexOf("freak")
$tr.Characters($idx + 1, 5).Text = "burn"

# 2) Re-assert "doesn't" as its own run (mirrors how PowerPoint
#    breaks runs around an edited word).
$text = $tr.Text
$idx = $text.IndexOf("doesn")
$tr.Characters($idx + 1, 7).Text = "doesn’t"

# 3) Split the space right after "burn" into its own run too.
$text = $tr.Text
$idx = $text.IndexOf("burn") + 4
$tr.Characters($idx + 1, 1).Text = " "
